$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3800
$ws.Range("I86").Value = 4083.3333
$ws.Range("J86").Value = 2950
$ws.Range("K86").Value = 4083.3333
$ws.Range("L86").Value = 2950
$ws.Range("M86").Value = -2960.3333
$ws.Range("N86").Value = -5196
$ws.Range("H89").Value = 3800
$ws.Range("I89").Value = 4083.3333
$ws.Range("J89").Value = 2950
$ws.Range("K89").Value = 20416.6665
$ws.Range("L89").Value = 14750
$ws.Range("M89").Value = -14800.6665
$ws.Range("N89").Value = -25982
$ws.Range("H138").Value = 2209.8484
$ws.Range("J138").Value = 2371.368
$ws.Range("L138").Value = 7114.103999999999
$ws.Range("N138").Value = -17394.104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1229
$ws.Range("I2").Value = 795.1818
$ws.Range("J2").Value = 2024.3334
$ws.Range("K2").Value = 795.1818
$ws.Range("L2").Value = 2024.3334
$ws.Range("M2").Value = -682.1818
$ws.Range("N2").Value = -2250.3334
$ws.Range("H15").Value = 3955
$ws.Range("J15").Value = 4899
$ws.Range("L15").Value = 4899
$ws.Range("N15").Value = -5599
$ws.Range("H32").Value = 10523.23
$ws.Range("J32").Value = 25796.916
$ws.Range("L32").Value = 25796.916
$ws.Range("N32").Value = -26370.916
$ws.Range("H61").Value = 62501308
$ws.Range("I61").Value = 100001040
$ws.Range("J61").Value = 1754.1666
$ws.Range("K61").Value = 100001040
$ws.Range("L61").Value = 1754.1666
$ws.Range("M61").Value = -100000828
$ws.Range("N61").Value = -2178.1666
$ws.Range("H74").Value = 2689.0908
$ws.Range("I74").Value = 2086.6667
$ws.Range("J74").Value = 5400
$ws.Range("K74").Value = 2086.6667
$ws.Range("L74").Value = 5400
$ws.Range("M74").Value = -1212.6667
$ws.Range("N74").Value = -7148
$ws.Range("H77").Value = 2689.0908
$ws.Range("I77").Value = 2086.6667
$ws.Range("J77").Value = 5400
$ws.Range("K77").Value = 10433.3335
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -6065.333500000001
$ws.Range("N77").Value = -35736
$ws.Range("H116").Value = 1229
$ws.Range("I116").Value = 795.1818
$ws.Range("J116").Value = 2024.3334
$ws.Range("K116").Value = 795.1818
$ws.Range("L116").Value = 2024.3334
$ws.Range("M116").Value = 1498.8182
$ws.Range("N116").Value = -6612.3334
$ws.Range("H136").Value = 62501308
$ws.Range("I136").Value = 100001040
$ws.Range("J136").Value = 1754.1666
$ws.Range("K136").Value = 300003120
$ws.Range("L136").Value = 5262.4998
$ws.Range("M136").Value = -300000570
$ws.Range("N136").Value = -10362.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1229
$ws.Range("I3").Value = 795.1818
$ws.Range("J3").Value = 2024.3334
$ws.Range("K3").Value = 795.1818
$ws.Range("L3").Value = 2024.3334
$ws.Range("M3").Value = -681.1818
$ws.Range("N3").Value = -2252.3334
$ws.Range("H123").Value = 50780
$ws.Range("J123").Value = 50780
$ws.Range("L123").Value = 50780
$ws.Range("N123").Value = -60580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4542.875
$ws.Range("I58").Value = 1265.2727
$ws.Range("K58").Value = 1265.2727
$ws.Range("M58").Value = -1062.2727
$ws.Range("H105").Value = 763.8
$ws.Range("I105").Value = 737.55554
$ws.Range("K105").Value = 737.55554
$ws.Range("M105").Value = 1009.44446
$ws.Range("H107").Value = 855.38464
$ws.Range("I107").Value = 342
$ws.Range("K107").Value = 342
$ws.Range("M107").Value = 1578
$ws.Range("H132").Value = 2521.12
$ws.Range("I132").Value = 2126.875
$ws.Range("J132").Value = 3222
$ws.Range("K132").Value = 6380.625
$ws.Range("L132").Value = 9666
$ws.Range("M132").Value = -3850.625
$ws.Range("N132").Value = -14726
$ws.Range("H136").Value = 4542.875
$ws.Range("I136").Value = 1265.2727
$ws.Range("K136").Value = 3795.8181
$ws.Range("M136").Value = -1245.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1125
$ws.Range("I17").Value = 1370
$ws.Range("K17").Value = 4110
$ws.Range("M17").Value = -3941
$ws.Range("H23").Value = 613
$ws.Range("I23").Value = 60
$ws.Range("J23").Value = 797.3333
$ws.Range("K23").Value = 180
$ws.Range("L23").Value = 2391.9999
$ws.Range("M23").Value = 55
$ws.Range("N23").Value = -2861.9999
$ws.Range("H55").Value = 2981.6365
$ws.Range("J55").Value = 3199.8
$ws.Range("L55").Value = 9599.400000000001
$ws.Range("N55").Value = -9953.400000000001
$ws.Range("H97").Value = 724.75
$ws.Range("I97").Value = 466.33334
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 1399.00002
$ws.Range("L97").Value = 4500
$ws.Range("M97").Value = -903.0000199999999
$ws.Range("N97").Value = -5492
$ws.Range("H106").Value = 2618.6
$ws.Range("I106").Value = 3035
$ws.Range("J106").Value = 2514.5
$ws.Range("K106").Value = 9105
$ws.Range("L106").Value = 7543.5
$ws.Range("M106").Value = -8159
$ws.Range("N106").Value = -9435.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 662.5714
$ws.Range("I97").Value = 565.5294
$ws.Range("J97").Value = 1075
$ws.Range("K97").Value = 565.5294
$ws.Range("L97").Value = 1075
$ws.Range("M97").Value = -69.52940000000001
$ws.Range("N97").Value = -2067
$ws.Range("H132").Value = 3171.0557
$ws.Range("I132").Value = 3009.4614
$ws.Range("J132").Value = 3591.2
$ws.Range("K132").Value = 9028.3842
$ws.Range("L132").Value = 10773.6
$ws.Range("M132").Value = -6498.3842
$ws.Range("N132").Value = -15833.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 22382.25
$ws.Range("J64").Value = 22382.25
$ws.Range("L64").Value = 22382.25
$ws.Range("N64").Value = -22832.25
$ws.Range("H67").Value = 22382.25
$ws.Range("J67").Value = 22382.25
$ws.Range("L67").Value = 22382.25
$ws.Range("N67").Value = -23942.25
$ws.Range("H122").Value = 20834782
$ws.Range("I122").Value = 27779310
$ws.Range("J122").Value = 1201.6666
$ws.Range("K122").Value = 83337930
$ws.Range("L122").Value = 3604.9998
$ws.Range("M122").Value = -83335480
$ws.Range("N122").Value = -8504.9998
$ws.Range("H123").Value = 40321.332
$ws.Range("J123").Value = 40321.332
$ws.Range("L123").Value = 40321.332
$ws.Range("N123").Value = -50121.332
$ws.Range("H136").Value = 1799.875
$ws.Range("I136").Value = 1119.8
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 3359.4
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -809.3999999999996
$ws.Range("N136").Value = -13899.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 15566.714
$ws.Range("J63").Value = 16790.166
$ws.Range("L63").Value = 16790.166
$ws.Range("N63").Value = -18038.166
$ws.Range("H66").Value = 15566.714
$ws.Range("J66").Value = 16790.166
$ws.Range("L66").Value = 50370.49800000001
$ws.Range("N66").Value = -56610.49800000001
$ws.Range("H113").Value = 460
$ws.Range("I113").Value = 218
$ws.Range("K113").Value = 654
$ws.Range("M113").Value = 1516
$ws.Range("H132").Value = 2927.2307
$ws.Range("I132").Value = 2656.875
$ws.Range("K132").Value = 7970.625
$ws.Range("M132").Value = -5440.625
$ws.Range("H136").Value = 1168.2609
$ws.Range("I136").Value = 931.0526
$ws.Range("K136").Value = 2793.1578
$ws.Range("M136").Value = -243.1578
